$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update the price in D29 (960 -> 2100)
$ws.Range("D29").Value = 2100
